$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$shp = $sm.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
Write-Output "Bold: $($tr.Font.Bold)"
$tr.Font.Bold = $tr.Font.Bold
